$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "95.330.81"
Set-TextValue "E2" "  -1.70%  "
Set-TextValue "D3" "3.612.58"
Set-TextValue "E3" "  -2.42%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "2.33"
Set-TextValue "E5" "  +21.27%  "
Set-TextValue "D6" "226.79"
Set-TextValue "E6" "  -4.54%  "
Set-TextValue "D7" "638.52"
Set-TextValue "E7" "  -2.65%  "
Set-TextValue "D8" "0.412"
Set-TextValue "E8" "  -3.44%  "
Set-TextValue "D9" "1.10"
Set-TextValue "E9" "  +2.44%  "
Set-TextValue "D10" "0.999"
Set-TextValue "E10" "  -0.03%  "
Set-TextValue "D11" "3.609.45"
Set-TextValue "E11" "  -2.47%  "
Set-TextValue "D12" "46.63"
Set-TextValue "E12" "  +5.30%  "
Set-TextValue "E13" "  -0.98%  "
Set-TextValue "D14" "0.0000291"
Set-TextValue "E14" "  -2.97%  "
Set-TextValue "D15" "6.48"
Set-TextValue "E15" "  -4.27%  "
Set-TextValue "D16" "4.283.67"
Set-TextValue "D17" "95.186.05"
Set-TextValue "E17" "  -1.75%  "
Set-TextValue "D18" "8.76"
Set-TextValue "E18" "  -2.00%  "
Set-TextValue "D19" "20.28"
Set-TextValue "E19" "  +8.78%  "
Set-TextValue "D20" "3.616.28"
Set-TextValue "E20" "  -2.28%  "
Set-TextValue "D21" "12.89"
Set-TextValue "E21" "  -0.93%  "
Set-TextValue "D22" "0.516"
Set-TextValue "E22" "  +1.29%  "
Set-TextValue "D23" "510.47"
Set-TextValue "E23" "  -2.52%  "
Set-TextValue "D24" "3.25"
Set-TextValue "E24" "  -5.19%  "
Set-TextValue "D25" "0.247"
Set-TextValue "E25" "  +27.70%  "
Set-TextValue "D26" "119.34"
Set-TextValue "E26" "  +17.11%  "
Set-TextValue "D27" "0.0000203"
Set-TextValue "E27" "  -3.72%  "
Set-TextValue "D28" "6.74"
Set-TextValue "E28" "  -2.63%  "
Set-TextValue "D29" "12.62"
Set-TextValue "E29" "  -5.97%  "
Set-TextValue "D30" "12.65"
Set-TextValue "E30" "  +2.04%  "
Set-TextValue "D31" "2.91"
Set-TextValue "E31" "  -3.19%  "
Set-TextValue "E32" "  -0.18%  "
Set-TextValue "E33" "  -0.02%  "
Set-TextValue "D36" "31.77"
Set-TextValue "E36" "  -1.63%  "
Set-TextValue "D37" "0.585"
Set-TextValue "E37" "  -2.44%  "
Set-TextValue "E38" "  -0.01%  "
Set-TextValue "D39" "596.63"
Set-TextValue "E39" "  -7.71%  "
Set-TextValue "D40" "8.31"
Set-TextValue "E40" "  -5.85%  "
Set-TextValue "D41" "6.77"
Set-TextValue "E41" "  -1.05%  "
Set-TextValue "D42" "40.77"
Set-TextValue "E42" "  +0.87%  "
Set-TextValue "E43" "  -1.64%  "
Set-TextValue "D44" "0.481"
Set-TextValue "E44" "  +6.60%  "
Set-TextValue "D45" "0.0475"
Set-TextValue "E45" "  +2.74%  "
Set-TextValue "E46" "  -7.21%  "
Set-TextValue "D47" "0.919"
Set-TextValue "E47" "  -3.90%  "
Set-TextValue "D48" "23.46"
Set-TextValue "E48" "  -0.73%  "
Set-TextValue "D49" "8.55"
Set-TextValue "E49" "  -0.26%  "
Set-TextValue "D50" "2.21"
Set-TextValue "E50" "  -3.33%  "

# Row 34: Fetch.AI -> Cronos
Set-TextValue "B34" "Cronos"
Set-TextValue "C34" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D34" "0.178"
Set-TextValue "E34" "  -5.18%  "

# Row 35: Cronos -> Fetch.AI
Set-TextValue "B35" "Fetch.AI"
Set-TextValue "C35" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D35" "1.79"
Set-TextValue "E35" "  -4.06%  "

# Row 51: OKB -> Aave
Set-TextValue "B51" "Aave"
Set-TextValue "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D51" "213.55"
Set-TextValue "E51" "  +4.18%  "
